# Update the "Förändrad" (Changed) date column C for all data rows
# from serial date 45181 (2023-09-12) to 45182 (2023-09-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows run from row 2 to row 385 in this workbook.
$firstRow = 2
$lastRow = 385

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45182
}
